$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.733.19"
$ws.Range("E2").Value = "  +1.61%  "

$ws.Range("D3").Value = "1.867.07"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'247.17"
$ws.Range("E5").Value = "  +2.21%  "

$ws.Range("D6").Value = "'0.7018"
$ws.Range("E6").Value = "  +0.83%  "

$ws.Range("D7").Value = "'0.9994"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.07768"
$ws.Range("E8").Value = "  -0.45%  "

$ws.Range("D9").Value = "'0.3093"
$ws.Range("E9").Value = "  +0.51%  "

$ws.Range("D10").Value = "'23.99"
$ws.Range("E10").Value = "  +0.34%  "

$ws.Range("D11").Value = "'0.07829"
$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").Value = "'5.183"
$ws.Range("E12").Value = "  +1.28%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.868.09"
$ws.Range("E13").Value = "  +0.39%  "

$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'93.24"
$ws.Range("E14").Value = "  +0.70%  "

$ws.Range("D15").Value = "'0.6979"
$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("D16").Value = "'6.662"
$ws.Range("E16").Value = "  +1.70%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "29.676.13"
$ws.Range("E17").Value = "  +1.45%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000008419"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").Value = "'243.97"
$ws.Range("E19").Value = "  -1.70%  "

$ws.Range("D20").Value = "2.112.19"
$ws.Range("E20").Value = "  -0.05%  "

$ws.Range("D21").Value = "'12.86"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").Value = "'0.9996"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").Value = "'7.611"
$ws.Range("E23").Value = "  +0.54%  "

$ws.Range("D24").Value = "'1.0000"
$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").Value = "'0.1533"
$ws.Range("E25").Value = "  +2.07%  "

$ws.Range("D26").Value = "'8.970"
$ws.Range("E26").Value = "  +0.90%  "

$ws.Range("D27").Value = "'160.20"
$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("D28").Value = "'18.46"
$ws.Range("E28").Value = "  -0.34%  "

$ws.Range("D29").Value = "'1.547"
$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("D30").Value = "'4.269"
$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("D31").Value = "'4.247"
$ws.Range("E31").Value = "  +0.67%  "

$ws.Range("D32").Value = "'1.197"
$ws.Range("E32").Value = "  -0.16%  "

$ws.Range("D33").Value = "'0.05158"
$ws.Range("E33").Value = "  -1.22%  "

$ws.Range("D34").Value = "'0.7910"
$ws.Range("E34").Value = "  +3.71%  "

$ws.Range("D35").Value = "'1.929"
$ws.Range("E35").Value = "  +4.16%  "

$ws.Range("D36").Value = "'1.160"
$ws.Range("E36").Value = "  -1.05%  "

$ws.Range("D37").Value = "'2.699"
$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("D38").Value = "1.342.60"
$ws.Range("E38").Value = "  +8.82%  "

$ws.Range("D39").Value = "'0.01883"
$ws.Range("E39").Value = "  +1.02%  "

$ws.Range("D40").Value = "'2.744"
$ws.Range("E40").Value = "  +0.54%  "

$ws.Range("D41").Value = "'0.9758"
$ws.Range("E41").Value = "  +7.12%  "

$ws.Range("D42").Value = "'6.114"
$ws.Range("E42").Value = "  +10.02%  "

$ws.Range("D43").Value = "'107.55"
$ws.Range("E43").Value = "  -1.93%  "

$ws.Range("D44").Value = "'0.9995"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").Value = "'9.822"
$ws.Range("E45").Value = "  +2.30%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'66.31"
$ws.Range("E46").Value = "  +1.65%  "

$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "2.009.12"
$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000124"
$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.5204"
$ws.Range("E49").Value = "  +0.44%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.785"
$ws.Range("E50").Value = "  +1.75%  "

$ws.Range("D51").Value = "'7.034"
$ws.Range("E51").Value = "  -0.15%  "
